$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge rows 2-6 (Sprouting Thrinax card) into a single A2 cell
$ws.Range("A2").Value = "('Sprouting Thrinax', ['{B}{R}{G}', 'Creature — Lizard', 'When Sprouting Thrinax dies, create three 1/1 green Saproling creature tokens.', '3/3'])"

# Merge rows 7-10 (Woolly Thoctar card) into a single A3 cell
$ws.Range("A3").Value = "('Woolly Thoctar', ['{R}{G}{W}', 'Creature — Beast', '5/4'])"

# Clear out the now-unused rows 4-10
$ws.Range("A4:A10").Clear()
